# Rename the single worksheet and restore the default (non-explicit) page
# setup scaling so the saved pageSetup element collapses back to just the
# orientation attribute (matches a vanilla Excel re-save).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Libro_IVA_Compra"
$ws.PageSetup.Zoom = $true
